$d = $word.ActiveDocument

# 1) First paragraph: "RAPPORT DE REUNION PARENTS-PROFS"
#    -> center alignment, bold paragraph mark, bold run text
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter
$p1.Range.Bold = 1
$p1.Range.Font.Bold = 1

# 2) Remove the "_GoBack" bookmark currently sitting right after
#    "Mr. BISIMWA prof de Français"
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3) Re-add the "_GoBack" bookmark right after "Shawn :" (collapsed,
#    i.e. zero-length, positioned before the paragraph mark).
#    This runtime mishandles Bookmarks.Add on a Range whose Start/End
#    sit exactly on a paragraph-mark boundary (it silently resets the
#    bookmark to the top of the document), so we work around it by
#    temporarily inserting a one-character marker at the target spot,
#    wrapping a non-collapsed bookmark around that marker (which works
#    reliably), and then deleting the marker text -- the bookmark
#    collapses down to the correct location in its place.
$findRange = $d.Content
$findRange.Find.Execute("Shawn :", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$insertPos = $findRange.End

$marker = $d.Range($insertPos, $insertPos)
$marker.InsertAfter([char]1)

$markerRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$markerRange2 = $d.Range($insertPos, $insertPos + 1)
$markerRange2.Text = ""
